$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: add P1=14, Q1=15 with the same style as O1 (bold + border) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update B2:I25 with new computed values ---
$ws.Range("B2").Value = 24.44307359118704
$ws.Range("C2").Value = 18.84945595138836
$ws.Range("D2").Value = 10.02066305241336
$ws.Range("E2").Value = 29.10526702065993
$ws.Range("F2").Value = 63.03962927701833
$ws.Range("G2").Value = 2.074527469938671
$ws.Range("H2").Value = 3.089966362669392
$ws.Range("I2").Value = 3.037267292020123

$ws.Range("B3").Value = 22.75403469952843
$ws.Range("C3").Value = 17.53033396951019
$ws.Range("D3").Value = 9.477274155190834
$ws.Range("E3").Value = 27.02775996347093
$ws.Range("F3").Value = 59.07048019432784
$ws.Range("G3").Value = 2.084514762941586
$ws.Range("H3").Value = 2.711583836633365
$ws.Range("I3").Value = 2.666704839883651

$ws.Range("B4").Value = 21.65952378285824
$ws.Range("C4").Value = 16.68690558818472
$ws.Range("D4").Value = 9.126601439550928
$ws.Range("E4").Value = 25.68821643124625
$ws.Range("F4").Value = 56.50830951348823
$ws.Range("G4").Value = 2.090765124095499
$ws.Range("H4").Value = 2.473648940128331
$ws.Range("I4").Value = 2.509508808433992

$ws.Range("B5").Value = 21.19839618377675
$ws.Range("C5").Value = 16.35017261237931
$ws.Range("D5").Value = 8.962248306382884
$ws.Range("E5").Value = 25.1250067760214
$ws.Range("F5").Value = 55.35697436933109
$ws.Range("G5").Value = 2.093373239817404
$ws.Range("H5").Value = 2.374520949779008
$ws.Range("I5").Value = 2.613255113535315

$ws.Range("B6").Value = 21.12034009970148
$ws.Range("C6").Value = 16.31276175098225
$ws.Range("D6").Value = 8.913715193426158
$ws.Range("E6").Value = 25.02928931015711
$ws.Range("F6").Value = 55.07122712887607
$ws.Range("G6").Value = 2.093843404499031
$ws.Range("H6").Value = 2.357148475066212
$ws.Range("I6").Value = 2.635588437125068

$ws.Range("B7").Value = 21.65169548959532
$ws.Range("C7").Value = 16.73292013436325
$ws.Range("D7").Value = 9.067926423508622
$ws.Range("E7").Value = 25.67736305216642
$ws.Range("F7").Value = 56.24481937174074
$ws.Range("G7").Value = 2.090896024981767
$ws.Range("H7").Value = 2.470075046960208
$ws.Range("I7").Value = 2.524493452909088

$ws.Range("B8").Value = 23.870143067299
$ws.Range("C8").Value = 18.46235242138411
$ws.Range("D8").Value = 9.765799130928215
$ws.Range("E8").Value = 28.39744764467337
$ws.Range("F8").Value = 61.38995501941691
$ws.Range("G8").Value = 2.078075537225312
$ws.Range("H8").Value = 2.957617360055809
$ws.Range("I8").Value = 2.908925381665187

$ws.Range("B9").Value = 27.77635781228673
$ws.Range("C9").Value = 21.51752086636726
$ws.Range("D9").Value = 11.09435833615976
$ws.Range("E9").Value = 33.24354161310295
$ws.Range("F9").Value = 70.89697344570428
$ws.Range("G9").Value = 2.053570185398751
$ws.Range("H9").Value = 3.877518438033591
$ws.Range("I9").Value = 3.818003966742009

$ws.Range("B10").Value = 30.36896984993525
$ws.Range("C10").Value = 23.60291155827736
$ws.Range("D10").Value = 11.73454473692673
$ws.Range("E10").Value = 35.62103610224544
$ws.Range("F10").Value = 76.20990935265392
$ws.Range("G10").Value = 2.036709757866656
$ws.Range("H10").Value = 4.482358605421413
$ws.Range("I10").Value = 4.449622519871673

$ws.Range("B11").Value = 31.36121202343197
$ws.Range("C11").Value = 24.30366210427836
$ws.Range("D11").Value = 9.970562816033523
$ws.Range("E11").Value = 29.02521488677462
$ws.Range("F11").Value = 69.78476296936714
$ws.Range("G11").Value = 2.03570407727512
$ws.Range("H11").Value = 4.790747882991368
$ws.Range("I11").Value = 4.510909130059678

$ws.Range("B12").Value = 31.68585877720369
$ws.Range("C12").Value = 24.44128593523203
$ws.Range("D12").Value = 8.470649040072542
$ws.Range("E12").Value = 23.0426619106675
$ws.Range("F12").Value = 63.59460163662632
$ws.Range("G12").Value = 2.037548408334782
$ws.Range("H12").Value = 5.512775492837027
$ws.Range("I12").Value = 4.456295562001463

$ws.Range("B13").Value = 31.52803904234584
$ws.Range("C13").Value = 24.22003956114179
$ws.Range("D13").Value = 7.009883965657012
$ws.Range("E13").Value = 17.05336036085701
$ws.Range("F13").Value = 56.73539812628097
$ws.Range("G13").Value = 2.041798004963339
$ws.Range("H13").Value = 6.450362237943956
$ws.Range("I13").Value = 4.308155851271346

$ws.Range("B14").Value = 31.20074615630609
$ws.Range("C14").Value = 23.91840789243268
$ws.Range("D14").Value = 6.025607134096393
$ws.Range("E14").Value = 12.91878041738623
$ws.Range("F14").Value = 51.49898527829699
$ws.Range("G14").Value = 2.045765259257151
$ws.Range("H14").Value = 7.201276089694248
$ws.Range("I14").Value = 4.16614545595082

$ws.Range("B15").Value = 31.01509405687257
$ws.Range("C15").Value = 23.77999269105106
$ws.Range("D15").Value = 5.778124340762592
$ws.Range("E15").Value = 11.90505057293832
$ws.Range("F15").Value = 49.98084492291337
$ws.Range("G15").Value = 2.047349627796285
$ws.Range("H15").Value = 7.37155626470092
$ws.Range("I15").Value = 4.109105430901421

$ws.Range("B16").Value = 29.98938299365462
$ws.Range("C16").Value = 23.00271204771266
$ws.Range("D16").Value = 5.728103560500506
$ws.Range("E16").Value = 11.57793924659197
$ws.Range("F16").Value = 48.56665406277664
$ws.Range("G16").Value = 2.053498711515744
$ws.Range("H16").Value = 7.055752869198356
$ws.Range("I16").Value = 3.877767251533557

$ws.Range("B17").Value = 29.36448092993919
$ws.Range("C17").Value = 22.55793443121987
$ws.Range("D17").Value = 6.19513540705
$ws.Range("E17").Value = 13.59355848671587
$ws.Range("F17").Value = 50.37157133274903
$ws.Range("G17").Value = 2.056209579269872
$ws.Range("H17").Value = 6.353811635908204
$ws.Range("I17").Value = 3.772518732282649

$ws.Range("B18").Value = 29.0326828153557
$ws.Range("C18").Value = 22.33313476896589
$ws.Range("D18").Value = 7.282784649532969
$ws.Range("E18").Value = 18.16992540045136
$ws.Range("F18").Value = 55.3707338003543
$ws.Range("G18").Value = 2.055887887705688
$ws.Range("H18").Value = 5.333293727927019
$ws.Range("I18").Value = 3.773529732633144

$ws.Range("B19").Value = 28.97345199471054
$ws.Range("C19").Value = 22.4033883500889
$ws.Range("D19").Value = 8.790375362995041
$ws.Range("E19").Value = 24.53208797344972
$ws.Range("F19").Value = 62.1806771869171
$ws.Range("G19").Value = 2.052752511625638
$ws.Range("H19").Value = 4.411514987705797
$ws.Range("I19").Value = 3.879578676368524

$ws.Range("B20").Value = 29.69820850587216
$ws.Range("C20").Value = 23.17917654727976
$ws.Range("D20").Value = 11.41151704201456
$ws.Range("E20").Value = 34.95363493126855
$ws.Range("F20").Value = 74.17869503357267
$ws.Range("G20").Value = 2.041462083656163
$ws.Range("H20").Value = 4.314225926561374
$ws.Range("I20").Value = 4.278525854025554

$ws.Range("B21").Value = 31.63458094084674
$ws.Range("C21").Value = 24.75784781552005
$ws.Range("D21").Value = 12.27927442675885
$ws.Range("E21").Value = 38.09980501401512
$ws.Range("F21").Value = 79.78419766314934
$ws.Range("G21").Value = 2.027249542232134
$ws.Range("H21").Value = 4.857352354950227
$ws.Range("I21").Value = 4.805519619189777

$ws.Range("B22").Value = 32.85232260185697
$ws.Range("C22").Value = 25.7093163221756
$ws.Range("D22").Value = 12.77195792536501
$ws.Range("E22").Value = 39.65220401541557
$ws.Range("F22").Value = 83.08140148289172
$ws.Range("G22").Value = 2.018208876805866
$ws.Range("H22").Value = 5.187829086509634
$ws.Range("I22").Value = 5.140100061090803

$ws.Range("B23").Value = 32.20893271165753
$ws.Range("C23").Value = 25.16305894123705
$ws.Range("D23").Value = 12.56410863771151
$ws.Range("E23").Value = 38.83252030482961
$ws.Range("F23").Value = 81.55397815487034
$ws.Range("G23").Value = 2.022918405131028
$ws.Range("H23").Value = 5.014325094539056
$ws.Range("I23").Value = 4.96347957896903

$ws.Range("B24").Value = 29.68185813993916
$ws.Range("C24").Value = 23.11238178875517
$ws.Range("D24").Value = 11.66976099974378
$ws.Range("E24").Value = 35.63223609375809
$ws.Range("F24").Value = 75.24053863113274
$ws.Range("G24").Value = 2.04088703992831
$ws.Range("H24").Value = 4.351378196554124
$ws.Range("I24").Value = 4.293550467567776

$ws.Range("B25").Value = 26.76275660477796
$ws.Range("C25").Value = 20.79094544270774
$ws.Range("D25").Value = 10.65831984704363
$ws.Range("E25").Value = 31.97679880725834
$ws.Range("F25").Value = 68.04404505137312
$ws.Range("G25").Value = 2.060308407478158
$ws.Range("H25").Value = 3.628166126820007
$ws.Range("I25").Value = 3.571692568677724

# --- Add P and Q columns (value 0, no special style) for rows 2-25 ---
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
